$d = $word.ActiveDocument

# Pull the whole package as WordprocessingML (pkg:package) text so we can
# perform precise, targeted edits to the style definitions that live in
# word/styles.xml. (The Styles collection's Font/ParagraphFormat proxies
# don't give us a way to fully remove a direct-formatting element such as
# <w:jc>, which is what this change requires for Author/Date.)
$xml = $d.WordOpenXML

# --- Title / TitleChar: drop the <w:spacing>/<w:kern> kerning pair from rPr ---
$xml = $xml.Replace(
    '<w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:spacing w:val="-10"/><w:kern w:val="28"/><w:sz w:val="56"/><w:szCs w:val="56"/>',
    '<w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:sz w:val="56"/><w:szCs w:val="56"/>')

# --- Author: base on Title, drop the now-inherited center alignment, and
#     set an explicit 12pt run size ---
$xml = $xml.Replace(
    '<w:style w:type="paragraph" w:customStyle="1" w:styleId="Author"><w:name w:val="Author"/><w:next w:val="BodyText"/><w:qFormat/><w:pPr><w:keepNext/><w:keepLines/><w:jc w:val="center"/></w:pPr></w:style>',
    '<w:style w:type="paragraph" w:customStyle="1" w:styleId="Author"><w:name w:val="Author"/><w:basedOn w:val="Title"/><w:next w:val="BodyText"/><w:qFormat/><w:pPr><w:keepNext/><w:keepLines/></w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:style>')

# --- Date: same treatment as Author ---
$xml = $xml.Replace(
    '<w:style w:type="paragraph" w:styleId="Date"><w:name w:val="Date"/><w:next w:val="BodyText"/><w:qFormat/><w:pPr><w:keepNext/><w:keepLines/><w:jc w:val="center"/></w:pPr></w:style>',
    '<w:style w:type="paragraph" w:styleId="Date"><w:name w:val="Date"/><w:basedOn w:val="Title"/><w:next w:val="BodyText"/><w:qFormat/><w:pPr><w:keepNext/><w:keepLines/></w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:style>')

$d.WordOpenXML = $xml
